$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for "Approximate size range (µm)" -- reuse the bold
# header formatting already used by A1:G1 (style index 3) via a
# format-only paste so we don't mint a duplicate font entry.
$ws.Range("I1").Value = "Approximate size range (µm)"
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data values for column I. Excel's text-autocomplete/number-guessing
# stamped these "N to M" strings with a date-ish display format even
# though they remain text values underneath.
$ws.Range("I2").Value = "10 to 85"
$ws.Range("I3").Value = "10 to 85"
$ws.Range("I4").Value = "10 to 85"
$ws.Range("I5").Value = "10 to 100 "
$ws.Range("I6").Value = "5 to 22"
$ws.Range("I7").Value = "6 to 10"
$ws.Range("I8").Value = "4 to 6 "
$ws.Range("I9").Value = "10 to 14"
$ws.Range("I10").Value = "12 to 14"
$ws.Range("I11").Value = "35 to 130"
$ws.Range("I12").Value = "18 to 26"
$ws.Range("I13").Value = "2 to 4"
$ws.Range("I14").Value = "6 to 8"
$ws.Range("I15").Value = "15 to 55"
$ws.Range("I16").Value = "4 to 15"
$ws.Range("I17").Value = "8 to 12"
$ws.Range("I18").Value = "3 to 5"
$ws.Range("I19").Value = "8 to 12"

# Rows 2-18 get "mmm-yy" (numFmtId 17); the last row (19) was typed
# slightly differently and picked up "d-mmm" (numFmtId 16).
$ws.Range("I2:I18").NumberFormat = "mmm-yy"
$ws.Range("I19").NumberFormat = "d-mmm"

$ws.Range("I19").Select() | Out-Null
